$wb = $excel.ActiveWorkbook

# --- Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")

# Version: 1.0.0 -> 1.0.1
$meta.Range("B3").Value = "1.0.1"

# Status: draft -> active
$meta.Range("B6").Value = "active"

# Experimental value cleared (row label "Experimental" stays, value removed)
$meta.Range("B7").ClearContents()

# Date: 2025-06-28 -> 2025-11-18
# (Assigning the literal string directly makes Excel auto-convert it into a date
#  serial number and also mints a brand-new cell style. To keep it as plain text
#  using the existing style, enter it as a text formula and then flatten the
#  formula down to a static value via copy/paste-values.)
$dateCell = $meta.Range("B8")
$dateCell.Formula = '="2025-11-18"'
$dateCell.Copy()
$dateCell.PasteSpecial(-4163)  # xlPasteValues

# --- Concepts sheet ---
$concepts = $wb.Worksheets.Item("Concepts")

# Add Definition text for MCS88126 (Pt—Bevidsthedsniveau...) row
$concepts.Range("D4").Value = "0: Habituel, 1: Agiteret eller reagerer kun på tale, 2: Reagerer kun på smerte, 3: Ingen reaktion"
